# added new data to PRVDRS
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Suicidio único
$ws.Range("D2").Value = 49
$ws.Range("I2").Value = 41

# Row 3 - Homicidio único
$ws.Range("H3").Value = 28
$ws.Range("I3").Value = 31

# Row 4 - Homicidio múltiple
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 13

# Row 10 - Total de víctimas mujeres: H10/I10 become computed totals like the rest of the row
$ws.Range("H10").Formula = "=SUM(H2:H8)"
$ws.Range("I10").Formula = "=SUM(I2:I8)"

# Update the saved cell selection
$ws.Range("I15").Select() | Out-Null
